# Insert a new data record as row 39, pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 39 (Excel shifts rows 39..129 down to 40..130,
# and the dimension / used range grows from A1:R129 to A1:R130 automatically).
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new record.
$ws.Range("A39").Value = 5
$ws.Range("B39").Value = "Macroferia Regional de Talca"
$ws.Range("C39").Value = "Maule"
$ws.Range("D39").Value = 45238
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = 300000000
$ws.Range("G39").Value = "Espárragos"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 4000
$ws.Range("K39").Value = 1300
$ws.Range("L39").Value = 1300
$ws.Range("M39").Value = 1300
$ws.Range("N39").Value = "$/kilo"
$ws.Range("O39").Value = "Provincia de Linares"
$ws.Range("P39").Value = 1300
$ws.Range("Q39").Value = 1
$ws.Range("R39").Value = "Hortaliza"
